# Update acquisition-datetime stamp (column A) on rows 2-11 of the
# "ランサーズ" sheet to reflect the latest scrape run.
# Commit message: Append: 2026-01-27 06:41 JST

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-27 06:41:37"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
